$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0.001
$ws.Range("K3").Value = 469
$ws.Range("L3").Value = 0.00469
